$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note above the table to reference the new source [5]
$ws.Range("A23").Value = "ayy, axy values for the ΔQy rms computation, using [5]"

# Add the new reference row (A42/B42), copying formatting from the row above
# (row 41, the previous "[4]" reference row) so the new row matches the
# existing reference-list styling.
$ws.Range("A41").Copy()
$ws.Range("A42").PasteSpecial(-4122)
$ws.Range("A42").Value = "[5]"

$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122)
$ws.Range("B42").Value = "https://github.com/natriant/CC_MD_2021/blob/master/helper_scripts/cmpt_tuneSpread_for_givenEmitandOctupoleSettings.py"

$ws.Hyperlinks.Add($ws.Range("B42"), "https://github.com/natriant/CC_MD_2021/blob/master/helper_scripts/cmpt_tuneSpread_for_givenEmitandOctupoleSettings.py")

# Adding the hyperlink re-applies Excel's built-in "Hyperlink" style; restore
# the reference-list formatting (matching B38/B39/B41) as the final step.
$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122)
